$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the corrected step counts for the affected trial rows.
# For these rows the y_corrSteps (E), y_nrSteps (G) and alienID (H)
# values were recalculated.
$ws.Range("E4").Value = 6
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 13

$ws.Range("E8").Value = 6
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 13

$ws.Range("E16").Value = 7
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 13

$ws.Range("E18").Value = 6
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 13

$ws.Range("E23").Value = 5
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 13

$ws.Range("E27").Value = 7
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 13

$ws.Range("A15").Select()
